$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 175
    3  = 59
    4  = 103
    5  = 73
    6  = 83
    7  = 91
    8  = 46
    9  = 119
    10 = 115
    11 = 231
    12 = 32
    13 = 112
    14 = 57
    15 = 125
    16 = 159
    17 = 40
    18 = 18
    19 = 113
    20 = 31
    21 = 15
    23 = 39
    24 = 23
    25 = 24
    26 = 27
    27 = 28
    28 = 25
    29 = 55
    30 = 67
    31 = 123
    32 = 94
    33 = 53
    34 = 43
    35 = 95
    36 = 7
    37 = 138
    38 = 85
    39 = 84
    40 = 79
    41 = 139
    42 = 70
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
